# CCC19 Derived Variables Spreadsheet
# 1. Insert a new "collapsed" cancer-treatment-timing variable row (T13a)
#    right after the existing T13 (cancer_tx_timing) row, shifting the
#    table / worksheet rows below it down by one.
# 2. Populate the new row's cells (Variable #, Variable Name, Category,
#    Description).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing T13 (cancer_tx_timing) row is worksheet row 43; push it (and
# everything below) down by one row to make room for the new T13a row.
$ws.Rows.Item(43).Insert()

$ws.Range("A43").Value = "T13a"
$ws.Range("B43").Value = "cancer_tx_timing_v2"
$ws.Range("C43").Value = "Cancer"
$ws.Range("D43").Value = "Timing of cancer treatment relative to COVID-19, collapsed"

# The dictionary lives inside an Excel Table (Table1); grow its range by
# one row so the new row is included (autoFilter / dimension follow).
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E135"))
